$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.729797666666667
$ws.Range("H2").Value = 5.189393000000001
$ws.Range("I2").Value = 0.06436583050179444
$ws.Range("J2").Value = 0.06436583050179444
$ws.Range("M2").Value = 2.294987
$ws.Range("N2").Value = 6.884961000000001
$ws.Range("O2").Value = 0.0158275801650097
$ws.Range("P2").Value = 0.0158275801650097
$ws.Range("Q2").Value = 3.969863157630335
$ws.Range("R2").Value = 35.72876841867301
$ws.Range("S2").Value = 0.001018755342154578
$ws.Range("T2").Value = 0.001018755342154578

$ws.Range("G3").Value = 1.729797666666667
$ws.Range("H3").Value = 5.189393000000001
$ws.Range("I3").Value = 0.06436583050179444
$ws.Range("J3").Value = 0.06436583050179444
$ws.Range("O3").Value = 0.769602070219672
$ws.Range("P3").Value = 0.7696020702196722
$ws.Range("Q3").Value = 193.0310807305418
$ws.Range("R3").Value = 1737.279726574876
$ws.Range("S3").Value = 0.0495360764055895
$ws.Range("T3").Value = 0.04953607640558952

$ws.Range("G4").Value = 1.729797666666667
$ws.Range("H4").Value = 5.189393000000001
$ws.Range("I4").Value = 0.06436583050179444
$ws.Range("J4").Value = 0.06436583050179444
$ws.Range("M4").Value = 31.11253633333333
$ws.Range("N4").Value = 93.337609
$ws.Range("O4").Value = 0.2145703496153182
$ws.Range("P4").Value = 0.2145703496153182
$ws.Range("Q4").Value = 53.8183927534819
$ws.Range("R4").Value = 484.3655347813371
$ws.Range("S4").Value = 0.01381099875405034
$ws.Range("T4").Value = 0.01381099875405035

$ws.Range("I5").Value = 0.2200595722726403
$ws.Range("J5").Value = 0.2200595722726403
$ws.Range("M5").Value = 2.294987
$ws.Range("N5").Value = 6.884961000000001
$ws.Range("O5").Value = 0.0158275801650097
$ws.Range("P5").Value = 0.0158275801650097
$ws.Range("Q5").Value = 13.57251792819934
$ws.Range("R5").Value = 122.152661353794
$ws.Range("S5").Value = 0.003483010521222961
$ws.Range("T5").Value = 0.003483010521222961

$ws.Range("I6").Value = 0.2200595722726403
$ws.Range("J6").Value = 0.2200595722726403
$ws.Range("O6").Value = 0.769602070219672
$ws.Range("P6").Value = 0.7696020702196722
$ws.Range("S6").Value = 0.1693583023926795
$ws.Range("T6").Value = 0.1693583023926795

$ws.Range("I7").Value = 0.2200595722726403
$ws.Range("J7").Value = 0.2200595722726403
$ws.Range("M7").Value = 31.11253633333333
$ws.Range("N7").Value = 93.337609
$ws.Range("O7").Value = 0.2145703496153182
$ws.Range("P7").Value = 0.2145703496153182
$ws.Range("S7").Value = 0.04721825935873782
$ws.Range("T7").Value = 0.04721825935873782

$ws.Range("H8").Value = 57.69206699999999
$ws.Range("I8").Value = 0.7155745972255653
$ws.Range("J8").Value = 0.7155745972255653
$ws.Range("M8").Value = 2.294987
$ws.Range("N8").Value = 6.884961000000001
$ws.Range("O8").Value = 0.0158275801650097
$ws.Range("P8").Value = 0.0158275801650097
$ws.Range("Q8").Value = 44.134181256043
$ws.Range("R8").Value = 397.207631304387
$ws.Range("S8").Value = 0.01132581430163216
$ws.Range("T8").Value = 0.01132581430163216

$ws.Range("H9").Value = 57.69206699999999
$ws.Range("I9").Value = 0.7155745972255653
$ws.Range("J9").Value = 0.7155745972255653
$ws.Range("O9").Value = 0.769602070219672
$ws.Range("P9").Value = 0.7696020702196722
$ws.Range("Q9").Value = 2145.985482808649
$ws.Range("R9").Value = 19313.86934527784
$ws.Range("S9").Value = 0.550707691421403
$ws.Range("T9").Value = 0.5507076914214031

$ws.Range("H10").Value = 57.69206699999999
$ws.Range("I10").Value = 0.7155745972255653
$ws.Range("J10").Value = 0.7155745972255653
$ws.Range("M10").Value = 31.11253633333333
$ws.Range("N10").Value = 93.337609
$ws.Range("O10").Value = 0.2145703496153182
$ws.Range("P10").Value = 0.2145703496153182
$ws.Range("Q10").Value = 598.3155102275337
$ws.Range("R10").Value = 5384.839592047802
$ws.Range("S10").Value = 0.15354109150253
$ws.Range("T10").Value = 0.1535410915025301
